$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: funding amount G3 changed 115 -> 67
$ws.Range("G3").Value = 67

# Row 4: D4 150 -> 138.35; remove F4/G4 (expense type + funding amount)
$ws.Range("D4").Value = 138.35
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""

# Row 5: D5 140 -> 160; remove F5/G5
$ws.Range("D5").Value = 160
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""

# New rows 8-11
$ws.Range("A8").Value = "Sonar Customs"
$ws.Range("B8").Value = "Iain P"
$ws.Range("D8").Value = 20

$ws.Range("A9").Value = "PCB Customs"
$ws.Range("D9").Value = 25

$ws.Range("A10").Value = "Paint"
$ws.Range("D10").Value = 20

$ws.Range("A11").Value = "Parts for PCB"
$ws.Range("D11").Value = 85.76

# Active cell selection moved to F18
$ws.Range("F18").Select()
